$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Add email and sms codes." -- the SMS/Email integration, template, and
# query rows (15, 16, 18, 19) were showing the stale "Service组件完成"
# status; bring them in line with the sibling rows (17, 20) that already
# read "Service完成" once those components' service layer was finished.
$ws.Range("D15").Value = "Service完成  总进度50%"
$ws.Range("D16").Value = "Service完成  总进度50%"
$ws.Range("D18").Value = "Service完成  总进度50%"
$ws.Range("D19").Value = "Service完成  总进度50%"

# Move the selection to where the work happened.
$ws.Range("D20").Select()
